$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data columns (F: 新能源汽车产销率, G: 新能源汽车销售量)
$ws.Range("F1:G17").Delete()

# Swap the "B" / "C" sub-rows within each year block (A:E only - the
# remaining columns after the F:G deletion above).
function Swap-Rows($r1, $r2) {
    for ($col = 1; $col -le 5; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $tmp = $c1.Value2
        $c1.Value2 = $c2.Value2
        $c2.Value2 = $tmp
    }
}

Swap-Rows 3 4
Swap-Rows 7 8
Swap-Rows 11 12
Swap-Rows 15 16

Write-Output "done"
